$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: update row 19 price value 15 -> 25 ---
$ws.Range("L19").Value = 25

# --- Step 2: insert a new row at 20 (shifts old row20->21, old row21->22) ---
$ws.Range("A19:N19").Copy()
$ws.Range("A20:N20").Insert()

# --- Step 3: row height for the new product row 20 ---
$ws.Rows.Item(20).RowHeight = 25.5

# --- Step 4: fill in the new row's content (item #17) ---
$ws.Range("A20").Value = 17
$ws.Range("B20").Value = "قصافات كبار"
$ws.Range("H20").Value = "9:0"
$ws.Range("L20").Value = 15
$ws.Range("N20").Value = 1
$ws.Range("C20:G20").Value = ""
$ws.Range("I20:K20").Value = ""
$ws.Range("M20").Value = ""

# --- Step 5: merges for row 20 (in case insert/copy dropped them) ---
if ($ws.Range("B20:G20").MergeCells -ne $true) { $ws.Range("B20:G20").Merge() }
if ($ws.Range("H20:K20").MergeCells -ne $true) { $ws.Range("H20:K20").Merge() }
if ($ws.Range("L20:M20").MergeCells -ne $true) { $ws.Range("L20:M20").Merge() }

# --- Step 6: the subtotal row (was row20, now row21) ---
$ws.Range("K21").Value = 671.5
$ws.Rows.Item(21).RowHeight = 25.5

Write-Host "edit complete"
